# Part 3 Sub Goals
# Sub goals for Part 3 of Problem Solving
#
# 1) In the "insight" sentence, remove the gramStart/gramEnd proofErr
#    markers around "is" and fold the text back into one continuous
#    sentence: "...because there is 5 fingers in a 10 count system. "
# 2) Add a new sentence after the "constraints" sentence of Part 3:
#    " The sub-goals are the different numbers in which she could
#    potentially count to. "

$d = $word.ActiveDocument

# --- Change 1: fix the "insight" sentence (drop proofErr wrapping "is") ---
$old1 = "The insight that I noticed while reading the problem is that the fingers would keep a consistent pattern of their numbering because there is 5 fingers in a 10 count system. "
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $old1, 2)

# --- Change 2: append the new sub-goals sentence after the constraints sentence ---
$old2 = "The constraints of this problem are that the counting system does not change, the pattern repeats over and over. "
$new2 = "The constraints of this problem are that the counting system does not change, the pattern repeats over and over.  The sub-goals are the different numbers in which she could potentially count to. "
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $new2, 2)

Write-Output "Change1 found: $found1"
Write-Output "Change2 found: $found2"
